$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '62.875.91'
$ws.Range('D3').Value = '2.465.38'
$ws.Range('E3').Value = '  +0.86%  '
Set-TextValue 'D5' '574.96'
$ws.Range('E5').Value = '  -0.24%  '
Set-TextValue 'D6' '147.72'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').Value = '2.464.61'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  +0.91%  '
Set-TextValue 'D14' '29.08'
$ws.Range('E14').Value = '  +2.51%  '
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').Value = '2.912.30'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '62.763.28'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '2.465.78'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('E19').Value = '  +0.30%  '
Set-TextValue 'D20' '11.02'
$ws.Range('E20').Value = '  -0.11%  '
Set-TextValue 'D21' '326.47'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('E22').Value = '  +8.25%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +0.00%  '
Set-TextValue 'D25' '10.06'
$ws.Range('E25').Value = '  +17.97%  '
Set-TextValue 'D26' '65.61'
$ws.Range('E26').Value = '  -1.14%  '
Set-TextValue 'D27' '644.95'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D28').Value = '0.0₃0984'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '2.582.91'
$ws.Range('E29').Value = '  +0.64%  '
Set-TextValue 'D30' '1.00'
$ws.Range('E30').Value = '  -14.96%  '
Set-TextValue 'D31' '1.44'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D38' '2.81'
$ws.Range('E38').Value = '  +3.35%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D39' '152.15'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D40' '0.369'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D41' '5.39'
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('B42').Value = 'EthereumClassic'
$ws.Range('C42').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D42' '18.70'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = '0.0₆0307'
$ws.Range('E44').Value = '  -36.07%  '
$ws.Range('E45').Value = '  -0.02%  '
Set-TextValue 'D46' '152.45'
$ws.Range('E46').Value = '  +4.98%  '
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('E48').Value = '  -1.25%  '
Set-TextValue 'D49' '20.52'
$ws.Range('E49').Value = '  -0.51%  '
Set-TextValue 'D50' '0.608'
$ws.Range('E50').Value = '  +0.52%  '
Set-TextValue 'D51' '0.0512'
$ws.Range('E51').Value = '  -0.79%  '
